$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "J2" = @("1924bced-b91c-472c-8c3c-30d51b7bef9e", "411bb346-23d4-4697-a4a5-d937017dc9b2")
    "J3" = @("f6bb4233-a6c3-495a-8e64-331b2883593a", "1030fceb-8cbc-46b8-b169-94fc19b82aee")
    "J4" = @("aa1f2b73-992f-4608-8b3a-2c44ceafe004", "3d789803-463a-405c-aef3-eade608ba9d9")
    "J5" = @("393b75be-0861-4ec9-b887-37269807287d", "93cd5286-281b-49b6-86f4-623b15162cc9")
}

foreach ($cellRef in $replacements.Keys) {
    $pair = $replacements[$cellRef]
    $oldUuid = $pair[0]
    $newUuid = $pair[1]
    $cell = $ws.Range($cellRef)
    $current = $cell.Text
    $updated = $current -replace [regex]::Escape($oldUuid), $newUuid
    $cell.Value = $updated
}
